# Rewrite the "confirmations" summary sheet from the old flat layout
# (generic "Summary" / "Confirmed" / "New nominations" / ... labels reused
# across every category) to the new layout where every label is spelled
# out per category (e.g. "     Civilian, New nominations"), and a new
# "Total new nominations" row is inserted right where the bare "Summary"
# divider used to be. The old trailing "Total withdrawn " row (42) is
# removed because its label/value moved into what is now row 41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused last row (old row 42, "Total withdrawn ") so the
# remaining rows shift up and the used range becomes A1:B41. Every row's
# existing number format travels with it untouched, which is what we
# want for every row except 37 and 40 (handled below).
$ws.Rows.Item(42).Delete()

# --- Column A labels (text) -------------------------------------------------
$labels = @{
    1  = "Labels"
    2  = "Congress"
    3  = "Session"
    4  = "Start Date"
    5  = "End Date"
    6  = "Civilian "
    7  = "     Civilian, New nominations"
    8  = "     Civilian, Carryover nominations"
    9  = "     Civilian, Confirmed "
    10 = "     Civilian, Unconfirmed "
    11 = "Civilian (FS, PHS, CG, NOAA)"
    12 = "     Civilian (FS, PHS, CG, NOAA), New nominations"
    13 = "     Civilian (FS, PHS, CG, NOAA), Carryover nominations"
    14 = "     Civilian (FS, PHS, CG, NOAA), Confirmed "
    15 = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed "
    16 = "Air Force "
    17 = "     Air Force, New nominations"
    18 = "     Air Force, Carryover nominations"
    19 = "     Air Force, Confirmed "
    20 = "     Air Force, Unconfirmed "
    21 = "Army "
    22 = "     Army, New nominations"
    23 = "     Army, Carryover nominations"
    24 = "     Army, Confirmed "
    25 = "     Army, Unconfirmed "
    26 = "Navy "
    27 = "     Navy, New nominations"
    28 = "     Navy, Carryover nominations"
    29 = "     Navy, Confirmed "
    30 = "     Navy, Unconfirmed "
    31 = "     Navy, Withdrawn "
    32 = "Marine Corps"
    33 = "     Marine Corps, New nominations"
    34 = "     Marine Corps, Carryover nominations"
    35 = "     Marine Corps, Confirmed "
    36 = "     Marine Corps, Unconfirmed "
    37 = "Total new nominations"
    38 = "Total carryover nominations"
    39 = "Total confirmed "
    40 = "Total unconfirmed       "
    41 = "Total withdrawn "
}

# --- Column B values (numbers / text, or $null for no value) ---------------
$values = @{
    1  = "Values"
    2  = 103
    3  = 2
    4  = 34359
    5  = 34669
    7  = 589
    8  = 172
    9  = 608
    10 = 153
    12 = 2440
    13 = 38
    14 = 2477
    15 = 1
    17 = 9124
    18 = 4
    19 = 9126
    20 = 2
    22 = 10278
    23 = 1686
    24 = 11960
    25 = 4
    27 = 11255
    28 = 657
    29 = 11901
    30 = 9
    31 = 2
    33 = 1359
    34 = 15
    35 = 1374
    36 = 0
    37 = 37446
    38 = 2572
    39 = 37446
    40 = 169
    41 = 2
}

# Row 37 used to be the bare "Summary" divider (no B value, General format).
# It now carries a number like the other "Total ..." rows, so copy B38's
# number format (thousands separator) onto it before writing the value.
$ws.Range("B38").Copy()
$ws.Range("B37").PasteSpecial(-4122)  # xlPasteFormats

# Row 40 used to be "Total confirmed " (thousands-separator format). It now
# holds "Total unconfirmed       ", which uses the plain integer format —
# copy B41's (still-correct) format onto it before overwriting the value.
$ws.Range("B41").Copy()
$ws.Range("B40").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

foreach ($r in 1..41) {
    $ws.Cells.Item($r, 1).Value = $labels[$r]

    if ($values.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $values[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = $null
    }
}
